# Final commit before submission
# - Mark issues 18-20 (rows 19-21) as completed ("y") with completion date 2019-09-22 (43730)
# - Update the saved cursor/selection position on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "completion" (E) column for rows 19-21
$ws.Range("E19").Value = "y"
$ws.Range("E20").Value = "y"
$ws.Range("E21").Value = "y"

# Copy the existing date cell's number format (style index reused, no new style
# created) onto the "date" (F) column for rows 19-21, then set the date value.
$ws.Range("F2").Copy()
$ws.Range("F19:F21").PasteSpecial(-4122)

$ws.Range("F19").Value = [DateTime]::FromOADate(43730)
$ws.Range("F20").Value = [DateTime]::FromOADate(43730)
$ws.Range("F21").Value = [DateTime]::FromOADate(43730)

# Update the sheet's stored selection/scroll position
[void]$ws.Range("I41").Select()
